# "New September Blog Section"
#
# 1) Refresh the cached "datetimeFigureOut" date-field text (7/22/20 -> 9/12/20)
#    on the slide master and every slide layout.
# 2) Move/resize the "Monthly Meetings (tbc)" textbox and reword it to
#    "Monthly Meetings (12th of each Month)" with a superscript "th".
# 3) Grow the "Meetings (tbc)" textbox and reword it to "Monthly Meetings (tbc)".

$p = $ppt.ActivePresentation

# --- 1) Date placeholder text: update on slide master + all custom layouts ---
$sm = $p.SlideMaster

for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "7/22/20") {
        $sh.TextFrame.TextRange.Text = "9/12/20"
    }
}

for ($L = 1; $L -le $sm.CustomLayouts.Count; $L++) {
    $cl = $sm.CustomLayouts.Item($L)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq "7/22/20") {
            $sh.TextFrame.TextRange.Text = "9/12/20"
        }
    }
}

# --- 2) & 3) Update the two meeting textboxes on slide 1 ---
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $t = $sh.TextFrame.TextRange.Text

        if ($t -eq "Monthly Meetings (tbc)") {
            # Reposition / resize
            $sh.Left = 255.43992614746094
            $sh.Top = 150.3428497314453
            $sh.Width = 70.49897766113281
            $sh.Height = 36.35157775878906

            # Reword with superscript "th"
            $sh.TextFrame.TextRange.Text = "Monthly Meetings (12th of each Month)"
            $tr = $sh.TextFrame.TextRange
            $sup = $tr.Characters(21, 2)
            $sup.Font.Superscript = $true
        }
        elseif ($t -eq "Meetings (tbc)") {
            $sh.Height = 36.35157775878906
            $sh.TextFrame.TextRange.Text = "Monthly Meetings (tbc)"
        }
    }
}
